$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column E a bit
$ws.Columns.Item(5).ColumnWidth = 14.625

# Add a new data row for the authorized person
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "12341234A"
$ws.Range("C2").Value = "Antonio"
$ws.Range("D2").Value = "González"
$ws.Range("E2").Value = "C/Almendra, 13"
$ws.Range("F2").NumberFormat = "mm-dd-yy"
$ws.Range("F2").Value = "10/10/1980"

# Update the selected cell
$ws.Range("H2").Select()
